$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 6403.2856
$ws.Range("I74").Value = 6564.6
$ws.Range("K74").Value = 6564.6
$ws.Range("M74").Value = -5628.6
$ws.Range("H77").Value = 6403.2856
$ws.Range("I77").Value = 6564.6
$ws.Range("K77").Value = 32823
$ws.Range("M77").Value = -28143
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 559539
$ws.Range("I32").Value = 4689.3057
$ws.Range("J32").Value = 1669238.5
$ws.Range("K32").Value = 4689.3057
$ws.Range("L32").Value = 1669238.5
$ws.Range("M32").Value = -4402.3057
$ws.Range("N32").Value = -1669812.5
$ws.Range("H61").Value = 11257.818
$ws.Range("I61").Value = 20307.2
$ws.Range("J61").Value = 3716.6667
$ws.Range("K61").Value = 20307.2
$ws.Range("L61").Value = 3716.6667
$ws.Range("M61").Value = -20095.2
$ws.Range("N61").Value = -4140.6667
$ws.Range("H122").Value = 63159840
$ws.Range("I122").Value = 109092376
$ws.Range("J122").Value = 2607.5
$ws.Range("K122").Value = 327277128
$ws.Range("L122").Value = 7822.5
$ws.Range("M122").Value = -327274678
$ws.Range("N122").Value = -12722.5
$ws.Range("H136").Value = 11257.818
$ws.Range("I136").Value = 20307.2
$ws.Range("J136").Value = 3716.6667
$ws.Range("K136").Value = 60921.60000000001
$ws.Range("L136").Value = 11150.0001
$ws.Range("M136").Value = -58371.60000000001
$ws.Range("N136").Value = -16250.0001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17545894
$ws.Range("I134").Value = 33335116
$ws.Range("J134").Value = 2314.2222
$ws.Range("K134").Value = 100005348
$ws.Range("L134").Value = 6942.6666
$ws.Range("M134").Value = -100002813
$ws.Range("N134").Value = -12012.6666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9047.423000000001
$ws.Range("I31").Value = 1064.2106
$ws.Range("J31").Value = 30716.143
$ws.Range("K31").Value = 1064.2106
$ws.Range("L31").Value = 30716.143
$ws.Range("M31").Value = -769.2106000000001
$ws.Range("N31").Value = -31306.143
$ws.Range("H34").Value = 9047.423000000001
$ws.Range("I34").Value = 1064.2106
$ws.Range("J34").Value = 30716.143
$ws.Range("K34").Value = 1064.2106
$ws.Range("L34").Value = 30716.143
$ws.Range("M34").Value = -862.2106000000001
$ws.Range("N34").Value = -31120.143
$ws.Range("H58").Value = 5723241.5
$ws.Range("I58").Value = 8404182
$ws.Range("J58").Value = 26243.5
$ws.Range("K58").Value = 8404182
$ws.Range("L58").Value = 26243.5
$ws.Range("M58").Value = -8403979
$ws.Range("N58").Value = -26649.5
$ws.Range("H62").Value = 2300
$ws.Range("I62").Value = 2200
$ws.Range("J62").Value = 2400
$ws.Range("K62").Value = 2200
$ws.Range("L62").Value = 2400
$ws.Range("M62").Value = -1576
$ws.Range("N62").Value = -3648
$ws.Range("H65").Value = 2300
$ws.Range("I65").Value = 2200
$ws.Range("J65").Value = 2400
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = -7880
$ws.Range("N65").Value = -18240
$ws.Range("H132").Value = 9808649
$ws.Range("I132").Value = 27778646
$ws.Range("J132").Value = 6832.3184
$ws.Range("K132").Value = 83335938
$ws.Range("L132").Value = 20496.9552
$ws.Range("M132").Value = -83333408
$ws.Range("N132").Value = -25556.9552
$ws.Range("H134").Value = 9767096
$ws.Range("I134").Value = 10418007
$ws.Range("J134").Value = 7814362.5
$ws.Range("K134").Value = 31254021
$ws.Range("L134").Value = 23443087.5
$ws.Range("M134").Value = -31251486
$ws.Range("N134").Value = -23448157.5
$ws.Range("H136").Value = 5723241.5
$ws.Range("I136").Value = 8404182
$ws.Range("J136").Value = 26243.5
$ws.Range("K136").Value = 25212546
$ws.Range("L136").Value = 78730.5
$ws.Range("M136").Value = -25209996
$ws.Range("N136").Value = -83830.5
$ws.Range("H138").Value = 35337.5
$ws.Range("J138").Value = 35337.5
$ws.Range("L138").Value = 35337.5
$ws.Range("N138").Value = -45617.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 200000880
$ws.Range("I118").Value = 250000600
$ws.Range("J118").Value = 2000
$ws.Range("K118").Value = 750001800
$ws.Range("L118").Value = 6000
$ws.Range("M118").Value = -750000557
$ws.Range("N118").Value = -8486
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H122").Value = 166668420
$ws.Range("I122").Value = 333334500
$ws.Range("K122").Value = 1000003500
$ws.Range("M122").Value = -1000001050
$ws.Range("H141").Value = 34157.25
$ws.Range("J141").Value = 34157.25
$ws.Range("L141").Value = 34157.25
$ws.Range("N141").Value = -44517.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6564.8237
$ws.Range("I46").Value = 690
$ws.Range("J46").Value = 14957.429
$ws.Range("K46").Value = 690
$ws.Range("L46").Value = 14957.429
$ws.Range("M46").Value = -502
$ws.Range("N46").Value = -15333.429
$ws.Range("H68").Value = 1610.2778
$ws.Range("I68").Value = 1502.4242
$ws.Range("J68").Value = 2796.6667
$ws.Range("K68").Value = 1502.4242
$ws.Range("L68").Value = 2796.6667
$ws.Range("M68").Value = -753.4241999999999
$ws.Range("N68").Value = -4294.6667
$ws.Range("H71").Value = 1610.2778
$ws.Range("I71").Value = 1502.4242
$ws.Range("J71").Value = 2796.6667
$ws.Range("K71").Value = 7512.120999999999
$ws.Range("L71").Value = 13983.3335
$ws.Range("M71").Value = -3768.120999999999
$ws.Range("N71").Value = -21471.3335
$ws.Range("H136").Value = 2795.724
$ws.Range("I136").Value = 2903.55
$ws.Range("K136").Value = 8710.650000000001
$ws.Range("M136").Value = -6160.650000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 57136344
$ws.Range("I132").Value = 60001120
$ws.Range("J132").Value = 51406796
$ws.Range("K132").Value = 180003360
$ws.Range("L132").Value = 154220388
$ws.Range("M132").Value = -180000830
$ws.Range("N132").Value = -154225448
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 21270.715
$ws.Range("J141").Value = 21270.715
$ws.Range("L141").Value = 21270.715
$ws.Range("N141").Value = -31630.715
